$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update nutritional / allergen breakdown values (sharedStrings content) ---

# Row 2 - Fries: diet no longer includes VEG
$ws.Range("E2").Value = "VGN,DF"

# Row 3 - Poutine: new ingredients description, diet, and nutrition label
$ws.Range("B3").Value = " Golden Crispy Fries / Vegetarian Brown Gravy / Cheese Curds"
$ws.Range("E3").Value = "VEG"
$ws.Range("F3").Value = "Poutine"

# Row 4 - Extra Cheese Curds: new ingredients, diet, and nutrition label
$ws.Range("B4").Value = "29% MF Milk / Enzymes / Salt"
$ws.Range("E4").Value = "GF,VEG"
$ws.Range("F4").Value = "Extra_Cheese_Curds"
$ws.Range("C4").WrapText = $true

# Row 5 - Side Gravy: nutrition label updated
$ws.Range("F5").Value = "Side_Gravy_-_Need_Update"

# --- Update selected cell to reflect the last-edited location ---
$ws.Range("F5").Select()
